# Applies the "update scripts wuth new tpm" edit to Rarres2-Cmklr1.xlsx:
#  - adds "Resolving-Mac" as a 4th sending cluster (rows 14-17)
#  - refreshes every ligand/receptor expression statistic (cols E:T) with new TPM-based numbers
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at row 14 to make room for the new "Resolving-Mac" sending-cluster block
$ws.Rows.Item(14).Resize(4).Insert()

$rowsData = @(
    @{ rn=2; A="ECs"; B="Rarres2"; C="Cmklr1"; D="ECs"; E=[double]"3"; F=[double]"1"; G=[double]"1.458525"; H=[double]"4.375575"; I=[double]"0.02525273220298681"; J=[double]"0.02525273220298681"; K=[double]"3"; L=[double]"1"; M=[double]"3.019277"; N=[double]"9.057831"; O=[double]"0.0408783237736836"; P=[double]"0.0408783237736836"; Q=[double]"4.403690986425"; R=[double]"39.63321887782499"; S=[double]"0.001032289363163821"; T=[double]"0.001032289363163821" },
    @{ rn=3; A="ECs"; B="Rarres2"; C="Cmklr1"; D="FAPs"; E=[double]"3"; F=[double]"1"; G=[double]"1.458525"; H=[double]"4.375575"; I=[double]"0.02525273220298681"; J=[double]"0.02525273220298681"; K=[double]"3"; L=[double]"1"; M=[double]"33.166404"; N=[double]"99.499212"; O=[double]"0.4490435959074953"; P=[double]"0.4490435959074953"; Q=[double]"48.37402939409999"; R=[double]"435.3662645469"; S=[double]"0.0113395776749182"; T=[double]"0.0113395776749182" },
    @{ rn=4; A="ECs"; B="Rarres2"; C="Cmklr1"; D="MuSCs"; E=[double]"3"; F=[double]"1"; G=[double]"1.458525"; H=[double]"4.375575"; I=[double]"0.02525273220298681"; J=[double]"0.02525273220298681"; K=[double]"3"; L=[double]"1"; M=[double]"1.677177666666666"; N=[double]"5.031533"; O=[double]"0.02270749311308342"; P=[double]"0.02270749311308342"; Q=[double]"2.446205556274999"; R=[double]"22.015850006475"; S=[double]"0.000573426242585863"; T=[double]"0.0005734262425858631" },
    @{ rn=5; A="ECs"; B="Rarres2"; C="Cmklr1"; D="Resolving-Mac"; E=[double]"3"; F=[double]"1"; G=[double]"1.458525"; H=[double]"4.375575"; I=[double]"0.02525273220298681"; J=[double]"0.02525273220298681"; K=[double]"3"; L=[double]"1"; M=[double]"35.99723933333333"; N=[double]"107.991718"; O=[double]"0.4873705872057377"; P=[double]"0.4873705872057378"; Q=[double]"52.50287349865"; R=[double]"472.52586148785"; S=[double]"0.01230743892231892"; T=[double]"0.01230743892231893" },
    @{ rn=6; A="FAPs"; B="Rarres2"; C="Cmklr1"; D="ECs"; E=[double]"3"; F=[double]"1"; G=[double]"48.53546666666667"; H=[double]"145.6064"; I=[double]"0.84033742450786"; J=[double]"0.8403374245078601"; K=[double]"3"; L=[double]"1"; M=[double]"3.019277"; N=[double]"9.057831"; O=[double]"0.0408783237736836"; P=[double]"0.0408783237736836"; Q=[double]"146.5420181909334"; R=[double]"1318.8781637184"; S=[double]"0.0343515853181757"; T=[double]"0.0343515853181757" },
    @{ rn=7; A="FAPs"; B="Rarres2"; C="Cmklr1"; D="FAPs"; E=[double]"3"; F=[double]"1"; G=[double]"48.53546666666667"; H=[double]"145.6064"; I=[double]"0.84033742450786"; J=[double]"0.8403374245078601"; K=[double]"3"; L=[double]"1"; M=[double]"33.166404"; N=[double]"99.499212"; O=[double]"0.4490435959074953"; P=[double]"0.4490435959074953"; Q=[double]"1609.7468957952"; R=[double]"14487.7220621568"; S=[double]"0.3773481388766528"; T=[double]"0.3773481388766529" },
    @{ rn=8; A="FAPs"; B="Rarres2"; C="Cmklr1"; D="MuSCs"; E=[double]"3"; F=[double]"1"; G=[double]"48.53546666666667"; H=[double]"145.6064"; I=[double]"0.84033742450786"; J=[double]"0.8403374245078601"; K=[double]"3"; L=[double]"1"; M=[double]"1.677177666666666"; N=[double]"5.031533"; O=[double]"0.02270749311308342"; P=[double]"0.02270749311308342"; Q=[double]"81.40260073457777"; R=[double]"732.6234066111999"; S=[double]"0.01908195627967849"; T=[double]"0.01908195627967849" },
    @{ rn=9; A="FAPs"; B="Rarres2"; C="Cmklr1"; D="Resolving-Mac"; E=[double]"3"; F=[double]"1"; G=[double]"48.53546666666667"; H=[double]"145.6064"; I=[double]"0.84033742450786"; J=[double]"0.8403374245078601"; K=[double]"3"; L=[double]"1"; M=[double]"35.99723933333333"; N=[double]"107.991718"; O=[double]"0.4873705872057377"; P=[double]"0.4873705872057378"; Q=[double]"1747.142809755022"; R=[double]"15724.2852877952"; S=[double]"0.409555744033353"; T=[double]"0.4095557440333531" },
    @{ rn=10; A="MuSCs"; B="Rarres2"; C="Cmklr1"; D="ECs"; E=[double]"3"; F=[double]"1"; G=[double]"7.712822"; H=[double]"23.138466"; I=[double]"0.1335389029981009"; J=[double]"0.133538902998101"; K=[double]"3"; L=[double]"1"; M=[double]"3.019277"; N=[double]"9.057831"; O=[double]"0.0408783237736836"; P=[double]"0.0408783237736836"; Q=[double]"23.287146069694"; R=[double]"209.584314627246"; S=[double]"0.005458846513138897"; T=[double]"0.005458846513138898" },
    @{ rn=11; A="MuSCs"; B="Rarres2"; C="Cmklr1"; D="FAPs"; E=[double]"3"; F=[double]"1"; G=[double]"7.712822"; H=[double]"23.138466"; I=[double]"0.1335389029981009"; J=[double]"0.133538902998101"; K=[double]"3"; L=[double]"1"; M=[double]"33.166404"; N=[double]"99.499212"; O=[double]"0.4490435959074953"; P=[double]"0.4490435959074953"; Q=[double]"255.806570432088"; R=[double]"2302.259133888792"; S=[double]"0.05996478919580944"; T=[double]"0.05996478919580946" },
    @{ rn=12; A="MuSCs"; B="Rarres2"; C="Cmklr1"; D="MuSCs"; E=[double]"3"; F=[double]"1"; G=[double]"7.712822"; H=[double]"23.138466"; I=[double]"0.1335389029981009"; J=[double]"0.133538902998101"; K=[double]"3"; L=[double]"1"; M=[double]"1.677177666666666"; N=[double]"5.031533"; O=[double]"0.02270749311308342"; P=[double]"0.02270749311308342"; Q=[double]"12.93577280537533"; R=[double]"116.421955248378"; S=[double]"0.003032333720158092"; T=[double]"0.003032333720158093" },
    @{ rn=13; A="MuSCs"; B="Rarres2"; C="Cmklr1"; D="Resolving-Mac"; E=[double]"3"; F=[double]"1"; G=[double]"7.712822"; H=[double]"23.138466"; I=[double]"0.1335389029981009"; J=[double]"0.133538902998101"; K=[double]"3"; L=[double]"1"; M=[double]"35.99723933333333"; N=[double]"107.991718"; O=[double]"0.4873705872057377"; P=[double]"0.4873705872057378"; Q=[double]"277.6402994693987"; R=[double]"2498.762695224588"; S=[double]"0.0650829335689945"; T=[double]"0.06508293356899451" },
    @{ rn=14; A="Resolving-Mac"; B="Rarres2"; C="Cmklr1"; D="ECs"; E=[double]"1"; F=[double]"0.3333333333333333"; G=[double]"0.05030299999999999"; H=[double]"0.150909"; I=[double]"0.0008709402910521559"; J=[double]"0.000870940291052156"; K=[double]"3"; L=[double]"1"; M=[double]"3.019277"; N=[double]"9.057831"; O=[double]"0.0408783237736836"; P=[double]"0.0408783237736836"; Q=[double]"0.151878690931"; R=[double]"1.366908218379"; S=[double]"3.560257920517625E-05"; T=[double]"3.560257920517626E-05" },
    @{ rn=15; A="Resolving-Mac"; B="Rarres2"; C="Cmklr1"; D="FAPs"; E=[double]"1"; F=[double]"0.3333333333333333"; G=[double]"0.05030299999999999"; H=[double]"0.150909"; I=[double]"0.0008709402910521559"; J=[double]"0.000870940291052156"; K=[double]"3"; L=[double]"1"; M=[double]"33.166404"; N=[double]"99.499212"; O=[double]"0.4490435959074953"; P=[double]"0.4490435959074953"; Q=[double]"1.668369620412"; R=[double]"15.015326583708"; S=[double]"0.0003910901601147806"; T=[double]"0.0003910901601147807" },
    @{ rn=16; A="Resolving-Mac"; B="Rarres2"; C="Cmklr1"; D="MuSCs"; E=[double]"1"; F=[double]"0.3333333333333333"; G=[double]"0.05030299999999999"; H=[double]"0.150909"; I=[double]"0.0008709402910521559"; J=[double]"0.000870940291052156"; K=[double]"3"; L=[double]"1"; M=[double]"1.677177666666666"; N=[double]"5.031533"; O=[double]"0.02270749311308342"; P=[double]"0.02270749311308342"; Q=[double]"0.08436706816633331"; R=[double]"0.7593036134969999"; S=[double]"1.97768706609737E-05"; T=[double]"1.977687066097371E-05" },
    @{ rn=17; A="Resolving-Mac"; B="Rarres2"; C="Cmklr1"; D="Resolving-Mac"; E=[double]"1"; F=[double]"0.3333333333333333"; G=[double]"0.05030299999999999"; H=[double]"0.150909"; I=[double]"0.0008709402910521559"; J=[double]"0.000870940291052156"; K=[double]"3"; L=[double]"1"; M=[double]"35.99723933333333"; N=[double]"107.991718"; O=[double]"0.4873705872057377"; P=[double]"0.4873705872057378"; Q=[double]"1.810769130184666"; R=[double]"16.296922171662"; S=[double]"0.0004244706810712253"; T=[double]"0.0004244706810712254" }
)

foreach ($row in $rowsData) {
    $ws.Cells.Item($row.rn, 1).Value = $row.A
    $ws.Cells.Item($row.rn, 2).Value = $row.B
    $ws.Cells.Item($row.rn, 3).Value = $row.C
    $ws.Cells.Item($row.rn, 4).Value = $row.D
    $ws.Cells.Item($row.rn, 5).Value = $row.E
    $ws.Cells.Item($row.rn, 6).Value = $row.F
    $ws.Cells.Item($row.rn, 7).Value = $row.G
    $ws.Cells.Item($row.rn, 8).Value = $row.H
    $ws.Cells.Item($row.rn, 9).Value = $row.I
    $ws.Cells.Item($row.rn, 10).Value = $row.J
    $ws.Cells.Item($row.rn, 11).Value = $row.K
    $ws.Cells.Item($row.rn, 12).Value = $row.L
    $ws.Cells.Item($row.rn, 13).Value = $row.M
    $ws.Cells.Item($row.rn, 14).Value = $row.N
    $ws.Cells.Item($row.rn, 15).Value = $row.O
    $ws.Cells.Item($row.rn, 16).Value = $row.P
    $ws.Cells.Item($row.rn, 17).Value = $row.Q
    $ws.Cells.Item($row.rn, 18).Value = $row.R
    $ws.Cells.Item($row.rn, 19).Value = $row.S
    $ws.Cells.Item($row.rn, 20).Value = $row.T
}
